# Update the cryptos list (Price / Volume(1h) columns, and one row re-ordering)
# as produced by the "Updated cryptos list" GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to keep its literal text (Excel would otherwise silently
    # re-interpret decimal-looking strings such as "324.99" as a number and
    # round-trip it as a float). We briefly mark the cell as Text, assign the
    # value, then restore the default "Normal" style so no extra formatting
    # lingers on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Map of row number -> (Price, Volume(1h)) new values.
# Rows 37 and 38 are handled separately below because the two coins
# (Hedera / InternetComputer(DFINITY)) swapped places.
$updates = @{
    2  = @("30.345.24", "  +1.25%  ")
    3  = @("2.009.83",  "  +4.84%  ")
    4  = @($null,       "  +0.12%  ")
    5  = @("324.99",    "  +1.48%  ")
    6  = @($null,       "  +0.09%  ")
    7  = @("0.5131",    "  +1.67%  ")
    8  = @("0.4267",    "  +5.84%  ")
    9  = @("0.08707",   "  +5.10%  ")
    10 = @($null,       "  +2.96%  ")
    11 = @("43.33",     "  +3.10%  ")
    12 = @("24.68",     "  +2.88%  ")
    13 = @("2.012.45",  "  +5.15%  ")
    14 = @("6.571",     "  +2.65%  ")
    15 = @("7.459",     "  +3.16%  ")
    16 = @("1.001",     "  -0.49%  ")
    17 = @("94.34",     "  +2.25%  ")
    19 = @("0.06531",   "  +0.47%  ")
    20 = @("18.89",     "  +3.89%  ")
    21 = @($null,       "  +0.05%  ")
    22 = @($null,       "  +4.39%  ")
    23 = @("30.413.39", $null)
    24 = @($null,       "  +4.53%  ")
    25 = @($null,       "  +2.81%  ")
    26 = @("2.248.15",  "  +5.26%  ")
    27 = @("22.39",     "  +0.89%  ")
    28 = @("162.23",    "  -0.15%  ")
    29 = @("2.425",     "  +4.98%  ")
    30 = @("131.03",    "  +1.51%  ")
    31 = @("1.138",     "  +0.32%  ")
    32 = @("0.1052",    "  +1.60%  ")
    33 = @($null,       "  +2.13%  ")
    34 = @("3.826",     $null)
    35 = @("1.373",     "  +14.88%  ")
    36 = @($null,       "  +3.31%  ")
    39 = @($null,       "  +8.76%  ")
    40 = @("9.158",     "  +4.92%  ")
    41 = @("0.2193",    "  +1.78%  ")
    42 = @("0.6644",    "  +2.95%  ")
    43 = @($null,       "  +2.24%  ")
    44 = @($null,       "  +0.06%  ")
    45 = @("13.63",     "  +2.02%  ")
    46 = @("0.6163",    "  +1.99%  ")
    47 = @("2.190",     "  -1.38%  ")
    49 = @($null,       "  +4.46%  ")
    50 = @("124.27",    "  +1.74%  ")
    51 = @("80.67",     "  +2.23%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $price = $vals[0]
    $volume = $vals[1]
    if ($null -ne $price) {
        Set-TextValue $ws.Range("D$row") $price
    }
    if ($null -ne $volume) {
        $ws.Range("E$row").Value = $volume
    }
}

# Rows 37 and 38 swapped: Hedera <-> InternetComputer(DFINITY)
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D37") "5.462"
$ws.Range("E37").Value = "  +0.71%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D38") "0.06674"
$ws.Range("E38").Value = "  +3.94%  "
